# The edit swaps the "Mедицина" and "Fisica" rows of the statistics table:
# what used to be row 2 (Mедицина) becomes row 3, and what used to be row 3
# (Fisica) becomes row 2 - all five columns (Profile, Avg. Score,
# Num. Students, Num. Universities, University names) move together.
#
# NOTE: reading back via ".Value" in this runtime yields a broken COM
# property-descriptor wrapper instead of the underlying scalar, so we read
# with ".Value2" (works for both numbers and strings) and only use the
# ".Value" setter for writing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = $ws.Range("A2:E2").Value2
$row3 = $ws.Range("A3:E3").Value2

$ws.Range("A2").Value = $row3[1, 1]
$ws.Range("B2").Value = $row3[1, 2]
$ws.Range("C2").Value = $row3[1, 3]
$ws.Range("D2").Value = $row3[1, 4]
$ws.Range("E2").Value = $row3[1, 5]

$ws.Range("A3").Value = $row2[1, 1]
$ws.Range("B3").Value = $row2[1, 2]
$ws.Range("C3").Value = $row2[1, 3]
$ws.Range("D3").Value = $row2[1, 4]
$ws.Range("E3").Value = $row2[1, 5]
